# Update the Price (D) and Volume(1h) (E) columns for rows 2-51 of the
# cryptos sheet with refreshed values. Numeric-looking Price values are
# prefixed with a literal apostrophe so Excel stores them as text (as in
# the original inline-string cells) instead of auto-converting to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    ,@(2, ('28.113.33'), '  +0.21%  ')
    ,@(3, ('1.867.56'), '  +3.75%  ')
    ,@(4, ("'" + '1.001'), '  -0.16%  ')
    ,@(5, ("'" + '311.41'), '  +0.87%  ')
    ,@(6, ("'" + '1.001'), '  -0.11%  ')
    ,@(7, ("'" + '0.5003'), '  -1.63%  ')
    ,@(8, ("'" + '0.3896'), '  +1.03%  ')
    ,@(9, ("'" + '0.09683'), '  +25.82%  ')
    ,@(10, ("'" + '1.134'), '  +3.46%  ')
    ,@(11, ("'" + '40.85'), '  +0.35%  ')
    ,@(12, ("'" + '6.453'), '  +1.93%  ')
    ,@(13, ("'" + '20.85'), '  +2.69%  ')
    ,@(14, ('1.867.73'), '  +3.79%  ')
    ,@(15, ("'" + '1.001'), '  -0.14%  ')
    ,@(16, ("'" + '7.375'), '  +1.64%  ')
    ,@(17, ("'" + '0.00001123'), '  +4.97%  ')
    ,@(18, ("'" + '93.01'), '  +1.00%  ')
    ,@(19, ("'" + '0.06591'), '  +0.32%  ')
    ,@(20, ("'" + '17.42'), '  +1.10%  ')
    ,@(21, ("'" + '1.001'), '  -0.14%  ')
    ,@(22, ("'" + '6.131'), '  +2.81%  ')
    ,@(23, ('28.175.03'), '  +0.33%  ')
    ,@(24, ("'" + '11.30'), '  +2.40%  ')
    ,@(25, ("'" + '2.279'), '  +1.76%  ')
    ,@(26, ("'" + '2.556'), '  +6.21%  ')
    ,@(27, ('2.078.78'), '  +3.38%  ')
    ,@(28, ("'" + '21.07'), '  +4.16%  ')
    ,@(29, ("'" + '157.44'), '  -1.50%  ')
    ,@(30, ("'" + '127.23'), '  +0.01%  ')
    ,@(31, ("'" + '0.1053'), '  -2.89%  ')
    ,@(32, ("'" + '1.058'), '  +1.30%  ')
    ,@(33, ("'" + '5.623'), '  +1.63%  ')
    ,@(34, ("'" + '3.625'), '  -0.65%  ')
    ,@(35, ("'" + '0.06742'), '  -3.24%  ')
    ,@(36, ("'" + '9.522'), '  +5.19%  ')
    ,@(37, ("'" + '0.02392'), '  +2.59%  ')
    ,@(38, ("'" + '0.2177'), '  +0.69%  ')
    ,@(39, ("'" + '11.47'), '  +0.23%  ')
    ,@(40, ("'" + '4.994'), '  -0.01%  ')
    ,@(41, ("'" + '0.6279'), '  +2.98%  ')
    ,@(42, ("'" + '1.170'), '  +1.96%  ')
    ,@(43, ("'" + '1.000'), '  -0.21%  ')
    ,@(44, ("'" + '13.50'), '  +2.06%  ')
    ,@(45, ("'" + '0.6013'), '  +2.23%  ')
    ,@(46, ("'" + '3.653'), '  -1.56%  ')
    ,@(47, ("'" + '1.257'), '  -3.24%  ')
    ,@(48, ("'" + '124.10'), '  -1.05%  ')
    ,@(49, ("'" + '1.974'), '  +2.82%  ')
    ,@(50, ("'" + '1.194'), '  +0.64%  ')
    ,@(51, ("'" + '0.06837'), '  +1.66%  ')
)

foreach ($entry in $updates) {
    $row = $entry[0]
    $priceText = $entry[1]
    $volumeText = $entry[2]
    $ws.Range("D$row").Value = $priceText
    $ws.Range("E$row").Value = $volumeText
}